{"js": "// Add a new \"Rule:\" bullet (top level) followed by a sub-bullet describing\n// the upload deadline, right after the last paragraph in the document\n// (\"Write process specifications for Process 1,3,4,6\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Build the two new list paragraphs as raw OOXML so the second paragraph's\n// text is split across three runs, exactly as authored, instead of being\n// collapsed into a single run by the higher level text-insertion APIs.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"0\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>Rule:</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"1\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>Any docume</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>nts must be uploaded before 11h45</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>pm.</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst insertionPoint = lastParagraph.getRange(Word.RangeLocation.after);\ninsertionPoint.insertOoxml(ooxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Add a new \"Rule:\" bullet (top level) followed by a sub-bullet describing\n# the upload deadline, right after the last paragraph in the document\n# (\"Write process specifications for Process 1,3,4,6\").\n\n$d = $word.ActiveDocument\n\n# Collapsed range at the very end of the document body (after the last\n# paragraph's mark). Building this from $d.Content.End (rather than\n# collapsing the last paragraph's own Range) keeps the existing text of\n# that last paragraph intact.\n$endPos = $d.Content.End\n$insertionRange = $d.Range($endPos, $endPos)\n\n# Insert the two new list paragraphs as raw OOXML so the second paragraph's\n# text lands in three separate runs, exactly as authored, instead of being\n# collapsed into a single run by plain text-insertion calls.\n$ooxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"0\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>Rule:</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"1\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>Any docume</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>nts must be uploaded before 11h45</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>pm.</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$insertionRange.InsertXML($ooxml)\n"}
